$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912" (schedule data refreshed, "Última actualización" bumped)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 03:59:53"

$ws1.Range("A6").Value = "03:59:53"
$ws1.Range("B6").Value = "04:01"
$ws1.Range("C6").Value = "81_EL PELIGRO"
$ws1.Range("D6").Value = 2

$ws1.Range("A7").Value = "03:59:53"
$ws1.Range("B7").Value = "04:46"
$ws1.Range("C7").Value = "215_EL PELIGRO"
$ws1.Range("D7").Value = 47

$ws1.Range("A8").Value = "03:59:53"
$ws1.Range("B8").Value = "04:53"
$ws1.Range("C8").Value = "11_ETCHEVERRY"
$ws1.Range("D8").Value = 54

$ws1.Range("A9").Value = "03:59:53"
$ws1.Range("B9").Value = "05:11"
$ws1.Range("C9").Value = "17_ROMERO"
$ws1.Range("D9").Value = 72

$ws1.Range("A10").Value = "03:59:53"
$ws1.Range("B10").Value = "05:21"
$ws1.Range("C10").Value = "23_HERNANDEZ"
$ws1.Range("D10").Value = 82

$ws1.Range("A11").Value = "03:59:53"
$ws1.Range("B11").Value = "05:31"
$ws1.Range("C11").Value = "81_EL PELIGRO"
$ws1.Range("D11").Value = 92

$ws1.Range("A12").Value = "03:59:53"
$ws1.Range("B12").Value = "05:47"
$ws1.Range("C12").Value = "14_ABASTO"
$ws1.Range("D12").Value = 108

$ws1.Range("A13").Value = "03:59:53"
$ws1.Range("B13").Value = "05:51"
$ws1.Range("C13").Value = "17_ROMERO"
$ws1.Range("D13").Value = 112

# ---------------------------------------------------------------------------
# Sheet "LP1912-215" (single row refreshed)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 03:59:53"

$ws2.Range("A6").Value = "03:59:53"
$ws2.Range("B6").Value = "04:46"
$ws2.Range("D6").Value = 47

# ---------------------------------------------------------------------------
# Sheet "6203-6173" (only timestamp refreshed, no data rows)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 03:59:53"
